# Updated cryptos list with GitHub Actions
# Refresh Price / Volume(1h) figures, and fix the swapped
# Arweave / Bittensor rows (44 & 45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.422.87'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').Value = '3.019.24'
$ws.Range('E3').Value = '  +2.16%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''595.83'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').Value = '''149.64'
$ws.Range('E6').Value = '  +4.83%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.018.06'
$ws.Range('E8').Value = '  +2.27%  '
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('D10').Value = '''6.40'
$ws.Range('E10').Value = '  +11.93%  '
$ws.Range('E11').Value = '  +3.92%  '
$ws.Range('D12').Value = '''0.459'
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('E13').Value = '  +3.43%  '
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('E15').Value = '  +2.64%  '
$ws.Range('D16').Value = '3.520.73'
$ws.Range('E16').Value = '  +1.88%  '
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').Value = '62.377.73'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('D19').Value = '2.996.50'
$ws.Range('E19').Value = '  +1.07%  '
$ws.Range('D20').Value = '''448.35'
$ws.Range('E20').Value = '  -0.31%  '
$ws.Range('E21').Value = '  +3.11%  '
$ws.Range('D22').Value = '''0.690'
$ws.Range('E22').Value = '  +1.85%  '
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('D24').Value = '''82.29'
$ws.Range('E24').Value = '  +1.39%  '
$ws.Range('E25').Value = '  +4.45%  '
$ws.Range('D26').Value = '''10.76'
$ws.Range('E26').Value = '  +12.85%  '
$ws.Range('E27').Value = '  -0.66%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E29').Value = '  +3.58%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').Value = '''7.16'
$ws.Range('E31').Value = '  +4.70%  '
$ws.Range('D32').Value = '''2.13'
$ws.Range('E32').Value = '  +3.63%  '
$ws.Range('D33').Value = '''27.57'
$ws.Range('E33').Value = '  +1.32%  '
$ws.Range('E34').Value = '  +2.75%  '
$ws.Range('D35').Value = '0.0₃0850'
$ws.Range('E35').Value = '  +8.90%  '
$ws.Range('E36').Value = '  +2.29%  '
$ws.Range('D37').Value = '''5.84'
$ws.Range('D38').Value = '''2.07'
$ws.Range('E38').Value = '  +0.20%  '
$ws.Range('E39').Value = '  +7.80%  '
$ws.Range('D40').Value = '''50.15'
$ws.Range('D41').Value = '''9.04'
$ws.Range('E41').Value = '  -0.46%  '
$ws.Range('D42').Value = '''0.122'
$ws.Range('E42').Value = '  +1.82%  '
$ws.Range('E43').Value = '  +7.92%  '
$ws.Range('B44').Value = 'Arweave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D44').Value = '''40.31'
$ws.Range('E44').Value = '  +10.11%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').Value = '''391.13'
$ws.Range('E45').Value = '  +0.49%  '
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').Value = '2.757.64'
$ws.Range('E47').Value = '  +1.52%  '
$ws.Range('D48').Value = '''133.73'
$ws.Range('E48').Value = '  +3.11%  '
$ws.Range('E50').Value = '  +1.11%  '
$ws.Range('E51').Value = '  -0.03%  '
